$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 103; this pushes the old row 103 (the footnote row)
# down to row 104 and copies formatting down from row 102 into the new row.
$ws.Rows.Item(103).Insert() | Out-Null

# Populate the newly inserted row 103 with the new day's data (2020-05-07).
$ws.Range("A103").Value = 43958
$ws.Range("B103").Value = 463
$ws.Range("C103").Value = 34703
$ws.Range("D103").Value = 77
$ws.Range("E103").Value = 7035

# Match the updated selection recorded on the sheet view.
$ws.Range("A103").Select() | Out-Null

# Update the print area defined name to cover the extra row (now through 105).
$wb.Names.Item(1).Delete() | Out-Null
$ws.Names.Add("_xlnm.Print_Area", '=相談件数!$A$1:$E$105') | Out-Null
